$d = $word.ActiveDocument

# --- Change 2 & 3: "... grinds the <m>inside of a large loaf of bread</m> ..."
#     becomes    "... grinds some <m>pith of coarse bread</m> ..."
# Locate the unique anchor phrase and compute sub-ranges within it so we can
# edit each run's text in place without disturbing the surrounding <m> tag
# runs (which must keep their own formatting/run boundaries).
$anchorText = "grinds the <m>inside of a large loaf of bread</m>"
$anchor = $d.Content
$anchor.Find.Execute($anchorText) | Out-Null
$base = $anchor.Start

$theOff = $anchorText.IndexOf(" the ")
$theLen = 5
$insideOff = $anchorText.IndexOf("inside")
$insideLen = 6
$ofText = " of a large loaf of bread"
$ofOff = $anchorText.IndexOf($ofText)
$ofLen = $ofText.Length

# Edit right-to-left so earlier (lower) offsets stay valid as text lengths change.
$ofRange = $d.Range($base + $ofOff, $base + $ofOff + $ofLen)
$ofRange.Text = "of coarse bread"

$insideRange = $d.Range($base + $insideOff, $base + $insideOff + $insideLen)
$insideRange.Text = "pith "

$theRange = $d.Range($base + $theOff, $base + $theOff + $theLen)
$theRange.Text = " some "

# --- Change 1: "<id>" + "p056v_1" + "</id>" (three runs) collapse into a
#     single run "<id>p056v_1</id>" (keeping the first run's formatting).
$idParagraph = $d.Paragraphs(6)
$idParagraph.Range.Find.Execute("<id>p056v_1</id>", $false, $false, $false, `
    $false, $false, $true, 1, $false, "<id>p056v_1</id>", 2) | Out-Null
